$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing row 74 with revised figures ---
$ws.Range("B74").Value = -26743
$ws.Range("C74").Value = 441966
$ws.Range("E74").Value = 197397
$ws.Range("G74").Value = 22167
$ws.Range("H74").Value = 96089
$ws.Range("I74").Value = 38651
$ws.Range("J74").Value = 224445
$ws.Range("K74").Value = 468709
$ws.Range("L74").Value = 33754
$ws.Range("M74").Value = 75224
$ws.Range("O74").Value = 52570
$ws.Range("Q74").Value = 19862
$ws.Range("R74").Value = 359731

# --- Add new row 75 with the latest quarter (01-04-2021) data ---
# Build the date label as a text formula first (string concat forces a text
# result instead of Excel auto-detecting "01-04-2021" as a date) and then
# convert it to a static value, so the cell ends up stored as plain text,
# matching the other "Serie" labels in column A.
$ws.Range("A75").Formula = '=""&"01-04-2021"'
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)

$ws.Range("B75").Value = -16701
$ws.Range("C75").Value = 454892
$ws.Range("D75").Value = 16477
$ws.Range("E75").Value = 206104
$ws.Range("F75").Value = 45303
$ws.Range("G75").Value = 25219
$ws.Range("H75").Value = 94556
$ws.Range("I75").Value = 41026
$ws.Range("J75").Value = 232312
$ws.Range("K75").Value = 471593
$ws.Range("L75").Value = 37860
$ws.Range("M75").Value = 72132
$ws.Range("N75").Value = 1265
$ws.Range("O75").Value = 49848
$ws.Range("P75").Value = 1258
$ws.Range("Q75").Value = 19760
$ws.Range("R75").Value = 361601
